# Commit: "Fruta / hortaliza, semanal"
# Inserts two new weekly price records at the top of the Betarraga series
# (new rows 250 and 251), shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 250; existing row 250 (and everything below)
# moves down to make room, preserving all of its data/formatting.
$ws.Rows.Item(250).Insert()
$ws.Rows.Item(250).Insert()

# New row 250
$ws.Range("A250").Value = 5
$ws.Range("B250").Value = "Macroferia Regional de Talca"
$ws.Range("C250").Value = "Maule"
$ws.Range("D250").Value = 44806
$ws.Range("E250").Value = 7
$ws.Range("F250").Value = 100114014
$ws.Range("G250").Value = "Betarraga"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 3000
$ws.Range("K250").Value = 1000
$ws.Range("L250").Value = 1000
$ws.Range("M250").Value = 1000
$ws.Range("N250").Value = "$/paquete 5 unidades"
$ws.Range("O250").Value = "Región del Maule"
$ws.Range("P250").Value = 200
$ws.Range("Q250").Value = 5
$ws.Range("R250").Value = "Hortaliza"

# New row 251
$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = 44806
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100114014
$ws.Range("G251").Value = "Betarraga"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Segunda"
$ws.Range("J251").Value = 3000
$ws.Range("K251").Value = 800
$ws.Range("L251").Value = 800
$ws.Range("M251").Value = 800
$ws.Range("N251").Value = "$/paquete 5 unidades"
$ws.Range("O251").Value = "Región del Maule"
$ws.Range("P251").Value = 160
$ws.Range("Q251").Value = 5
$ws.Range("R251").Value = "Hortaliza"
